$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update totals / counters
$ws.Range("E11").Value = 171954
$ws.Range("F13").Value = 3

# 2) Swap the two existing period codes (2507 / 2506 -> 2506 / 2507) and update
#    the "Salario Basico" amounts for both existing rows.
$ws.Range("E16").Value = "2506"
$ws.Range("G16").Value = 1432964
$ws.Range("E17").Value = "2507"
$ws.Range("G17").Value = 1432964

# 3) Insert a new blank row before row 18 - this pushes the trailing
#    signature block (rows 18-23) down by one row. Row 17 (still carrying
#    the old "closing" / bottom-border formatting) keeps its position.
$ws.Rows.Item(18).Insert()

# 4) Grab that "closing" formatting from row 17 before it gets overwritten,
#    and stamp it onto the brand new row 18.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# ...then restyle row 17 itself with the normal "middle of table" formatting
# (matching row 16), since it is no longer the last data row.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# 5) Fill in the new account-statement row with the new period (2508).
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143364121"
$ws.Range("D18").Value = "LEONELA MARTINEZ HERRERA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 57318
$ws.Range("G18").Value = 1432964

$excel.CutCopyMode = 0
